$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 3275541.5
$ws.Range("I43").Value = 12954.556
$ws.Range("J43").Value = 6945952
$ws.Range("K43").Value = 12954.556
$ws.Range("L43").Value = 6945952
$ws.Range("M43").Value = -12885.556
$ws.Range("N43").Value = -6946090

$ws.Range("H110").Value = 31418.334
$ws.Range("J110").Value = 31418.334
$ws.Range("L110").Value = 31418.334
$ws.Range("N110").Value = -39598.334

$ws.Range("H111").Value = 3039.3076
$ws.Range("I111").Value = 2565.125
$ws.Range("J111").Value = 3798
$ws.Range("K111").Value = 7695.375
$ws.Range("L111").Value = 11394
$ws.Range("M111").Value = -4628.375
$ws.Range("N111").Value = -17528

$ws.Range("H113").Value = 2709.4375
$ws.Range("I113").Value = 2496
$ws.Range("K113").Value = 2496
$ws.Range("M113").Value = 758

$ws.Range("H137").Value = 1016.2
$ws.Range("I137").Value = 842.2222
$ws.Range("J137").Value = 1211.925
$ws.Range("K137").Value = 2526.6666
$ws.Range("L137").Value = 3635.775
$ws.Range("M137").Value = 23.33339999999998
$ws.Range("N137").Value = -8735.775

$ws.Range("H138").Value = 558667.75
$ws.Range("I138").Value = 852.6486
$ws.Range("J138").Value = 1062061.9
$ws.Range("K138").Value = 2557.9458
$ws.Range("L138").Value = 3186185.7
$ws.Range("M138").Value = 2582.0542
$ws.Range("N138").Value = -3196465.7

$ws.Range("H141").Value = 544.6818
$ws.Range("I141").Value = 544.6818
$ws.Range("K141").Value = 1634.0454
$ws.Range("M141").Value = 3545.9546

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4163.157
$ws.Range("I32").Value = 3765.3064
$ws.Range("K32").Value = 3765.3064
$ws.Range("M32").Value = -3478.3064

$ws.Range("H61").Value = 21277526
$ws.Range("I61").Value = 23256638
$ws.Range("J61").Value = 2078.5
$ws.Range("K61").Value = 23256638
$ws.Range("L61").Value = 2078.5
$ws.Range("M61").Value = -23256426
$ws.Range("N61").Value = -2502.5

$ws.Range("H74").Value = 953.1951
$ws.Range("I74").Value = 696.30554
$ws.Range("K74").Value = 696.30554
$ws.Range("M74").Value = 177.69446

$ws.Range("H77").Value = 953.1951
$ws.Range("I77").Value = 696.30554
$ws.Range("K77").Value = 3481.5277
$ws.Range("M77").Value = 886.4723000000004

$ws.Range("H102").Value = 83334824
$ws.Range("I102").Value = 166666670
$ws.Range("J102").Value = 2981
$ws.Range("K102").Value = 166666670
$ws.Range("L102").Value = 2981
$ws.Range("M102").Value = -166665048
$ws.Range("N102").Value = -6225

$ws.Range("H132").Value = 2757.2354
$ws.Range("I132").Value = 2818.8696
$ws.Range("J132").Value = 2628.3635
$ws.Range("K132").Value = 8456.6088
$ws.Range("L132").Value = 7885.0905
$ws.Range("M132").Value = -5926.6088
$ws.Range("N132").Value = -12945.0905

$ws.Range("H136").Value = 21277526
$ws.Range("I136").Value = 23256638
$ws.Range("J136").Value = 2078.5
$ws.Range("K136").Value = 69769914
$ws.Range("L136").Value = 6235.5
$ws.Range("M136").Value = -69767364
$ws.Range("N136").Value = -11335.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 5367.811
$ws.Range("I134").Value = 2022.037
$ws.Range("J134").Value = 14401.4
$ws.Range("K134").Value = 6066.111
$ws.Range("L134").Value = 43204.2
$ws.Range("M134").Value = -3531.111
$ws.Range("N134").Value = -48274.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1967.0454
$ws.Range("I31").Value = 2246.6875
$ws.Range("J31").Value = 1221.3334
$ws.Range("K31").Value = 2246.6875
$ws.Range("L31").Value = 1221.3334
$ws.Range("M31").Value = -1951.6875
$ws.Range("N31").Value = -1811.3334

$ws.Range("H34").Value = 1967.0454
$ws.Range("I34").Value = 2246.6875
$ws.Range("J34").Value = 1221.3334
$ws.Range("K34").Value = 2246.6875
$ws.Range("L34").Value = 1221.3334
$ws.Range("M34").Value = -2044.6875
$ws.Range("N34").Value = -1625.3334

$ws.Range("H93").Value = 30000
$ws.Range("J93").Value = 30000
$ws.Range("L93").Value = 30000
$ws.Range("N93").Value = -33744

$ws.Range("H132").Value = 3219.64
$ws.Range("I132").Value = 3216.7334
$ws.Range("J132").Value = 3224
$ws.Range("K132").Value = 9650.200199999999
$ws.Range("L132").Value = 9672
$ws.Range("M132").Value = -7120.200199999999
$ws.Range("N132").Value = -14732

$ws.Range("H134").Value = 26317016
$ws.Range("I134").Value = 1268.7693
$ws.Range("J134").Value = 83334470
$ws.Range("K134").Value = 3806.3079
$ws.Range("L134").Value = 250003410
$ws.Range("M134").Value = -1271.3079
$ws.Range("N134").Value = -250008480

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 9257.477000000001
$ws.Range("I3").Value = 4395.8335
$ws.Range("J3").Value = 15739.667
$ws.Range("K3").Value = 13187.5005
$ws.Range("L3").Value = 47219.001
$ws.Range("M3").Value = -13075.5005
$ws.Range("N3").Value = -47443.001

$ws.Range("H113").Value = 720.10345
$ws.Range("I113").Value = 543
$ws.Range("J113").Value = 726.4286
$ws.Range("K113").Value = 1629
$ws.Range("L113").Value = 2179.2858
$ws.Range("M113").Value = 541
$ws.Range("N113").Value = -6519.2858

$ws.Range("H140").Value = 21778.02
$ws.Range("I140").Value = 57266.39
$ws.Range("J140").Value = 2990.0588
$ws.Range("K140").Value = 171799.17
$ws.Range("L140").Value = 8970.1764
$ws.Range("M140").Value = -166619.17
$ws.Range("N140").Value = -19330.1764

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 573.74194
$ws.Range("I107").Value = 813.5333000000001
$ws.Range("J107").Value = 348.9375
$ws.Range("K107").Value = 813.5333000000001
$ws.Range("L107").Value = 348.9375
$ws.Range("M107").Value = 1106.4667
$ws.Range("N107").Value = -4188.9375

$ws.Range("H132").Value = 1951.1034
$ws.Range("I132").Value = 1548.55
$ws.Range("J132").Value = 2845.6667
$ws.Range("K132").Value = 4645.65
$ws.Range("L132").Value = 8537.000100000001
$ws.Range("M132").Value = -2115.65
$ws.Range("N132").Value = -13597.0001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 5326.6924
$ws.Range("I46").Value = 5424.5
$ws.Range("J46").Value = 5283.222
$ws.Range("K46").Value = 5424.5
$ws.Range("L46").Value = 5283.222
$ws.Range("M46").Value = -5236.5
$ws.Range("N46").Value = -5659.222

$ws.Range("H132").Value = 18666.29
$ws.Range("I132").Value = 1205.9143
$ws.Range("J132").Value = 44129.332
$ws.Range("K132").Value = 3617.7429
$ws.Range("L132").Value = 132387.996
$ws.Range("M132").Value = -1087.7429
$ws.Range("N132").Value = -137447.996

$ws.Range("H135").Value = 36415.89
$ws.Range("J135").Value = 36415.89
$ws.Range("L135").Value = 36415.89
$ws.Range("N135").Value = -46555.89

$ws.Range("H136").Value = 1116.9062
$ws.Range("I136").Value = 1025.5518
$ws.Range("J136").Value = 2000
$ws.Range("K136").Value = 3076.6554
$ws.Range("L136").Value = 6000
$ws.Range("M136").Value = -526.6553999999996
$ws.Range("N136").Value = -11100

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 523.1667
$ws.Range("I107").Value = 488.66666
$ws.Range("J107").Value = 626.6667
$ws.Range("K107").Value = 1465.99998
$ws.Range("L107").Value = 1880.0001
$ws.Range("M107").Value = 454.0000199999999
$ws.Range("N107").Value = -5720.0001

$ws.Range("H132").Value = 2513.628
$ws.Range("I132").Value = 2569.6572
$ws.Range("K132").Value = 7708.971600000001
$ws.Range("M132").Value = -5178.971600000001

$ws.Range("H136").Value = 567.9783
$ws.Range("I136").Value = 408.8857
$ws.Range("J136").Value = 1074.1818
$ws.Range("K136").Value = 1226.6571
$ws.Range("L136").Value = 3222.5454
$ws.Range("M136").Value = 1323.3429
$ws.Range("N136").Value = -8322.545399999999
